$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# NitroXBots (sheet "NitroXBots") gains a new Test Case 24 row: a bunch of
# new "bot" related header columns (G..AC) are inserted conceptually, the
# previously-trailing columns (Side/Quantity/OpenOrderNumber/SkipAtStepNum/
# SIT) move out to the new tail (AD..AH), and the old Position/Price/Random/
# OrderType/DealRef columns are dropped. Row 2 gets populated with the new
# QA_TestCase_Auto_NitroX_024 data.
# ---------------------------------------------------------------------------

$ws = $wb.Worksheets.Item("NitroXBots")

# -- Row 1 headers -----------------------------------------------------------
# Columns A-F keep their existing header text/style; copy header formatting
# (style 9 = centered + wrapped + bordered) onto the new header cells G1:AH1.
$ws.Range("A1").Copy()
$ws.Range("G1:AH1").PasteSpecial(-4122)

$ws.Range("G1").Value  = "Service"
$ws.Range("H1").Value  = "Method"
$ws.Range("I1").Value  = "Bot Quantity"
$ws.Range("J1").Value  = "Instrument Type"
$ws.Range("K1").Value  = "Order Direction"
$ws.Range("L1").Value  = "Min Time Break"
$ws.Range("M1").Value  = "Max Time Break"
$ws.Range("N1").Value  = "Random Range"
$ws.Range("O1").Value  = "Order Amount"
$ws.Range("P1").Value  = "Execution Strategy"
$ws.Range("Q1").Value  = "Price Increment"
$ws.Range("R1").Value  = "Reserved Amount"
$ws.Range("S1").Value  = "Min Price"
$ws.Range("T1").Value  = "Max Price"
$ws.Range("U1").Value  = "Updating Break"
$ws.Range("V1").Value  = "Depth Level"
$ws.Range("W1").Value  = "Spread Benchmark"
$ws.Range("X1").Value  = "Target Altcoin"
$ws.Range("Y1").Value  = "Target Quotecoin"
$ws.Range("Z1").Value  = "Trigger Condition"
$ws.Range("AA1").Value = "Order Type"
$ws.Range("AB1").Value = "Stop Condition"
$ws.Range("AC1").Value = "Deal Ref"
$ws.Range("AD1").Value = "Side"
$ws.Range("AE1").Value = "Quantity"
$ws.Range("AF1").Value = "OpenOrderNumber"
$ws.Range("AG1").Value = "SkipAtStepNum"
$ws.Range("AH1").Value = "SIT"

# -- Row 2 data ---------------------------------------------------------------
# Copy the data-row formatting (style 6 = wrapped + bordered) from an
# existing populated cell onto the whole new data range F2:AH2.
$ws.Range("A2").Copy()
$ws.Range("F2:AH2").PasteSpecial(-4122)

$ws.Range("A2").Value  = "QA_TestCase_Auto_NitroX_024"
$ws.Range("B2").Value  = "Spot"
$ws.Range("C2").Value  = "Trader01@Tinyex"
$ws.Range("D2").Value  = "ETH"
$ws.Range("E2").Value  = "USDT"
$ws.Range("G2").Value  = "apl_bot"
$ws.Range("H2").Value  = "execution_bot"
$ws.Range("I2").Value  = 1
$ws.Range("K2").Value  = "BUY"
$ws.Range("L2").Value  = 15
$ws.Range("M2").Value  = 45
$ws.Range("O2").Value  = 10
$ws.Range("P2").Value  = "VANILLA"
$ws.Range("S2").Value  = 3000
$ws.Range("T2").Value  = 6000
$ws.Range("AA2").Value = "LIMIT"
$ws.Range("AH2").Value = $false

# ---------------------------------------------------------------------------
# View state: selection moves to NitroXBots, which becomes the active tab.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("NitroXBuySell")
$ws.Activate()
$ws.Range("T6").Select()
